# Apply the "actualice las tareas" commit to _TAREAS.xlsx
#
# Summary of the change (see xml diff):
#  - Clear the comment text in D5 and D6 on the "Tareas diarias" sheet
#    (the two related shared strings disappear from sharedStrings.xml
#    once nothing references them any more).
#  - Row 5 and Row 6 go back to their default (non custom) height since
#    their long comments are gone.
#  - Row 14 keeps a custom height, but it shrinks from 30 to 19.5.
#  - The sheet view no longer scrolls to show row 4 at the top, and the
#    active selection moves from D7 to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tareas diarias")

# --- Clear the two obsolete comments -------------------------------------
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()

# --- Row heights -----------------------------------------------------------
# Rows 5 and 6 no longer need the extra height that the removed comments
# required, so let them fall back to the sheet's standard height.
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()

# Row 14's comment got shorter too, so its custom height shrinks.
$ws.Rows.Item(14).RowHeight = 19.5

# --- Sheet view / selection --------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D5").Select()
